# Excel COM-interop script: rename "RACINE" -> "ROUX" throughout the
# "Astreintes 2024" / "Astreintes 2025" sheets (standalone surname cells,
# "Jean RACINE" mentions, and the "Jean RACINE" inside longer comment
# texts), then restore the expected active sheet / selection state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Astreintes 2024")
$ws2 = $wb.Worksheets.Item("Astreintes 2025")

# Replace every occurrence of RACINE with ROUX on both sheets. This covers
# the standalone "RACINE" surname cells, "Jean RACINE" cells, and the
# occurrences embedded inside longer comment sentences.
$ws1.Cells.Replace("RACINE", "ROUX")
$ws2.Cells.Replace("RACINE", "ROUX")

# The workbook now opens on "Astreintes 2025" with F22 selected, and
# "Astreintes 2024" should become the active sheet with I12 selected.
$ws2.Activate()
$ws2.Range("F22").Select()

$ws1.Activate()
$ws1.Range("I12").Select()
